# Apply the rebrand/version-bump edit described by the diff:
#  - IBM -> LinuxForHealth rebrand of the FHIR URLs
#  - Version bump 7.0.0 -> 8.0.0
#  - Date update
#  - Publisher rename
#  - Clear stale Constraint(s) text on the root "Extension" element row

$wb = $excel.ActiveWorkbook

# --- "Metadata" sheet (Property / Value overview) ---
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/attributed-provider-with-period"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- "Elements" sheet (per-element structure definition table) ---
$elements = $wb.Worksheets.Item("Elements")

# Extension.url row (row 5): Fixed Value column (Q)
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/attributed-provider-with-period"

# Extension.value[x] row (row 6): Type(s) column (J)
$elements.Range("J6").Value = "Reference {http://linuxforhealth.org/fhir/cdm/StructureDefinition/reference-with-period}
"

# Extension.value[x].extension / referencePeriod row (row 9): Type(s) column (J)
$elements.Range("J9").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/reference-period}
"

# Root "Extension" row (row 2): clear the stale Constraint(s) text (column AI)
$elements.Range("AI2").Value = ""
